# Update the cryptos price list with freshly scraped values (GitHub Actions run).
# Values that look numeric (e.g. "405.00", "0.0000126") are prefixed with a
# leading apostrophe so Excel stores them as text, matching the source data
# which keeps prices as literal strings (including multi-dot formats like
# "60.875.13").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.875.13"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.361.60"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'405.00"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "'134.98"
$ws.Range("E6").Value = "  +11.25%  "
$ws.Range("D7").Value = "'0.589"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.672"
$ws.Range("E9").Value = "  +5.54%  "
$ws.Range("D10").Value = "'0.119"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "'42.37"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "3.886.57"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "'8.29"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "'19.51"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "3.354.39"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "60.874.13"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'10.93"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'0.0000126"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").Value = "'83.50"
$ws.Range("E22").Value = "  +9.59%  "
$ws.Range("D23").Value = "'307.37"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").Value = "'12.63"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'4.77"
$ws.Range("E26").Value = "  +11.81%  "
$ws.Range("D27").Value = "'8.30"
$ws.Range("E27").Value = "  +9.52%  "
$ws.Range("D28").Value = "'29.31"
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("D29").Value = "'7.44"
$ws.Range("E29").Value = "  -7.49%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'11.24"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "'41.10"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("D35").Value = "'2.49"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "'0.0478"
$ws.Range("D37").Value = "'51.85"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'3.40"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'137.03"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.97"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "'4.00"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").Value = "'0.288"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").Value = "'16.56"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "'21.31"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "2.116.41"
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("E50").Value = "  -4.32%  "
$ws.Range("D51").Value = "'1.89"
$ws.Range("E51").Value = "  -0.55%  "
